# Add the new Modules-directory / README timecard entries for 2/12/2020 and 2/13/2020.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 2020-02-12 (serial 43873) block: rows 24-25 ---------------------------
# Copy the formatting (date number format, etc.) from an existing date cell
# so the new date cell reuses the same cell style as the others in column A.
$ws.Range("A19").Copy($ws.Range("A24"))
$ws.Range("A24").Value = 43873
$ws.Range("B24").Value = "Create script to parallelize writing to Zarr files (Viz2)"
$ws.Range("B25").Value = "Use dask distributed rather than multirprocessing to accomplish this task"

# --- 2020-02-13 (serial 43874) block: rows 27-29 ----------------------------
$ws.Range("A19").Copy($ws.Range("A27"))
$ws.Range("A27").Value = 43874
$ws.Range("B27").Value = "Curate scripts for parallel writing to Zarr files"
$ws.Range("B28").Value = "Create simple script to read Zarr files to either Zarr intermediate or dask"
$ws.Range("B29").Value = "Curate repository and update README"

# --- Column A width (best-fit to fit the date strings) ---------------------
$ws.Columns.Item(1).ColumnWidth = 9.7109375

# --- View state: scroll so row 3 is at the top, select B29 -----------------
$win = $wb.Windows.Item(1)
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("B29").Select()
